$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.62
$ws.Range("J2").Value = 3.9
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.74
$ws.Range("S2").Value = 2.9
$ws.Range("T2").Value = 1.87
$ws.Range("W2").Value = 2.6
$ws.Range("AB2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AN2").Value = 8.199999999999999

# Row 3 updates
$ws.Range("F3").Value = 8.199999999999999
$ws.Range("I3").Value = 1.46
$ws.Range("M3").Value = 1.03
$ws.Range("P3").Value = 2.36
$ws.Range("R3").Value = 1.54
$ws.Range("U3").Value = 1.97
$ws.Range("V3").Value = 3.15
$ws.Range("AB3").Value = 1000
$ws.Range("AG3").Value = 1000

# Row 4 updates
$ws.Range("J4").Value = 1.2
$ws.Range("O4").Value = 1.2
$ws.Range("S4").Value = 1.19
$ws.Range("W4").Value = 1.16
